$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '74.744.38'
$ws.Range('E2').Value = '  +8.92%  '
$ws.Range('D3').Value = '2.596.80'
$ws.Range('E3').Value = '  +6.95%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '186.35'
$ws.Range('E5').Value = '  +16.20%  '
$ws.Range('D6').Value = '582.18'
$ws.Range('E6').Value = '  +4.18%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').Value = '0.534'
$ws.Range('E8').Value = '  +5.00%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').Value = '0.207'
$ws.Range('E9').Value = '  +24.99%  '
$ws.Range('D10').Value = '2.592.02'
$ws.Range('E10').Value = '  +6.73%  '
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('D12').Value = '0.361'
$ws.Range('E12').Value = '  +9.24%  '
$ws.Range('D13').Value = '4.79'
$ws.Range('E13').Value = '  +4.18%  '
$ws.Range('E14').Value = '  +10.33%  '
$ws.Range('D15').Value = '74.628.37'
$ws.Range('E15').Value = '  +8.91%  '
$ws.Range('D16').Value = '3.045.47'
$ws.Range('E16').Value = '  +5.85%  '
$ws.Range('D17').Value = '26.29'
$ws.Range('E17').Value = '  +13.76%  '
$ws.Range('D18').Value = '2.557.81'
$ws.Range('E18').Value = '  +5.18%  '
$ws.Range('D19').Value = '8.99'
$ws.Range('E19').Value = '  +30.10%  '
$ws.Range('D20').Value = '11.76'
$ws.Range('E20').Value = '  +12.29%  '
$ws.Range('D21').Value = '377.22'
$ws.Range('E21').Value = '  +12.63%  '
$ws.Range('D22').Value = '2.31'
$ws.Range('E22').Value = '  +19.82%  '
$ws.Range('E23').Value = '  +6.48%  '
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('D25').Value = '69.83'
$ws.Range('E25').Value = '  +4.45%  '
$ws.Range('D26').Value = '4.19'
$ws.Range('E26').Value = '  +14.36%  '
$ws.Range('D27').Value = '9.30'
$ws.Range('E27').Value = '  +13.53%  '
$ws.Range('D28').Value = '2.730.08'
$ws.Range('E28').Value = '  +6.78%  '
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -1.06%  '
$ws.Range('D30').Value = '0.0₃0949'
$ws.Range('E30').Value = '  +16.06%  '
$ws.Range('D31').Value = '511.35'
$ws.Range('E31').Value = '  +19.97%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').Value = '1.39'
$ws.Range('E32').Value = '  +21.27%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = '7.94'
$ws.Range('E33').Value = '  +11.14%  '
$ws.Range('D34').Value = '1.74'
$ws.Range('E34').Value = '  +7.80%  '
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('E36').Value = '  +13.11%  '
$ws.Range('D37').Value = '158.85'
$ws.Range('E37').Value = '  -1.13%  '
$ws.Range('D38').Value = '19.20'
$ws.Range('E38').Value = '  +7.15%  '
$ws.Range('D39').Value = '19.36'
$ws.Range('E39').Value = '  +1.91%  '
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('E41').Value = '  +13.73%  '
$ws.Range('E42').Value = '  +12.82%  '
$ws.Range('E43').Value = '  +8.24%  '
$ws.Range('E44').Value = '  +20.61%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '155.44'
$ws.Range('E45').Value = '  +18.48%  '
$ws.Range('B46').Value = 'ImmutableX'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D46').Value = '1.17'
$ws.Range('E46').Value = '  +9.20%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').Value = '38.93'
$ws.Range('E47').Value = '  +4.30%  '
$ws.Range('D48').Value = '0.0815'
$ws.Range('E48').Value = '  +13.80%  '
$ws.Range('E49').Value = '  +8.64%  '
$ws.Range('D50').Value = '0.523'
$ws.Range('E50').Value = '  +8.91%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '19.98'
$ws.Range('E51').Value = '  +18.70%  '
